$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 97-98; this shifts the existing rows 97..226 down to 99..228
# (Excel copies formatting from the row above, which already matches the D-column date style)
$ws.Rows("97:98").Insert()

# Row 97: new "Primera" quality entry for the latest week (date 44571)
$ws.Range("A97").Value2 = 1
$ws.Range("B97").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C97").Value2 = "Arica y Parinacota"
$ws.Range("D97").Value2 = 44571
$ws.Range("E97").Value2 = 15
$ws.Range("F97").Value2 = 100114014
$ws.Range("G97").Value2 = "Betarraga"
$ws.Range("H97").Value2 = "Sin especificar"
$ws.Range("I97").Value2 = "Primera"
$ws.Range("J97").Value2 = 1200
$ws.Range("K97").Value2 = 300
$ws.Range("L97").Value2 = 350
$ws.Range("M97").Value2 = 325
$ws.Range("N97").Value2 = "$/paquete 4 unidades"
$ws.Range("O97").Value2 = "Región de Arica y Parinacota"
$ws.Range("P97").Value2 = 81
$ws.Range("Q97").Value2 = 4
$ws.Range("R97").Value2 = "Hortaliza"

# Row 98: new "Segunda" quality entry for the latest week (date 44571)
$ws.Range("A98").Value2 = 1
$ws.Range("B98").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C98").Value2 = "Arica y Parinacota"
$ws.Range("D98").Value2 = 44571
$ws.Range("E98").Value2 = 15
$ws.Range("F98").Value2 = 100114014
$ws.Range("G98").Value2 = "Betarraga"
$ws.Range("H98").Value2 = "Sin especificar"
$ws.Range("I98").Value2 = "Segunda"
$ws.Range("J98").Value2 = 1200
$ws.Range("K98").Value2 = 300
$ws.Range("L98").Value2 = 350
$ws.Range("M98").Value2 = 325
$ws.Range("N98").Value2 = "$/paquete 5 unidades"
$ws.Range("O98").Value2 = "Región de Arica y Parinacota"
$ws.Range("P98").Value2 = 65
$ws.Range("Q98").Value2 = 5
$ws.Range("R98").Value2 = "Hortaliza"
